# The deck's single slide master/theme (ppt/theme/theme1.xml, "Integral" /
# "Red Violet") gets its 12 theme colours replaced with the stock Office
# theme palette ("Office Theme" / "Office"):
#   dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink
#
# RGB() isn't a built-in here (that's VBA), so build the packed 0x00BBGGRR
# COM colour value ourselves.
function ToComRgb($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Target "Office Theme" colour scheme, in clrScheme order.
$officeColors = @(
    @(0, 0, 0),        # 1  dk1      -> 000000
    @(255, 255, 255),  # 2  lt1      -> FFFFFF
    @(68, 84, 106),    # 3  dk2      -> 44546A
    @(231, 230, 230),  # 4  lt2      -> E7E6E6
    @(91, 155, 213),   # 5  accent1  -> 5B9BD5
    @(237, 125, 49),   # 6  accent2  -> ED7D31
    @(165, 165, 165),  # 7  accent3  -> A5A5A5
    @(255, 192, 0),    # 8  accent4  -> FFC000
    @(68, 114, 196),   # 9  accent5  -> 4472C4
    @(112, 173, 71),   # 10 accent6  -> 70AD47
    @(5, 99, 193),     # 11 hlink    -> 0563C1
    @(149, 79, 114)    # 12 folHlink -> 954F72
)

$p = $ppt.ActivePresentation

# Go through the theme colour scheme exposed on a slide (backed by the same
# ppt/theme/theme1.xml the slide master uses) and overwrite every slot.
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $rgb = $officeColors[$i - 1]
    $themeColors.Colors($i).RGB = ToComRgb $rgb[0] $rgb[1] $rgb[2]
}
